$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 26, shifting the existing rows 26:97 down to 27:98
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new record.
# Columns that stay constant across every data row in this sheet:
$ws.Range("A26").Value = 11
$ws.Range("B26").Value = "Vega Monumental Concepción"
$ws.Range("C26").Value = "Bíobío"
$ws.Range("E26").Value = 8
$ws.Range("F26").Value = 100112021
$ws.Range("G26").Value = "Ají"
$ws.Range("I26").Value = "Primera"
$ws.Range("R26").Value = "Hortaliza"

# Columns with the new record's actual data.
$ws.Range("D26").Value = 44672
$ws.Range("H26").Value = "Chilena(o)"
$ws.Range("J26").Value = 35
$ws.Range("K26").Value = 24000
$ws.Range("L26").Value = 25000
$ws.Range("M26").Value = 24571
$ws.Range("N26").Value = "$/saco 25 kilos"
$ws.Range("O26").Value = "Región Metropolitana"
$ws.Range("P26").Value = 983
$ws.Range("Q26").Value = 25
